$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings (e.g. "1.005") that must stay as plain text
# and not be auto-converted to numbers by Excel smart-entry parsing. Force
# the cell to Text format before assigning, then restore the default style
# so no stray number-format style is left attached to the cell.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.869.63"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.707.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.49%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.27%  "

$ws.Range("E6").Value = "  +0.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3732"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.54%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "49.38"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.25%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3432"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.53%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.219"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.50%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07537"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.60%  "

$ws.Range("E12").Value = "  +0.18%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.22"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.28%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.308"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.76%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.101"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.11%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.707.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.46%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001131"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.88%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06724"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.68%  "

$ws.Range("E19").Value = "  +0.02%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "83.99"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.50%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.28"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.81%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.376"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.11%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "13.08"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +7.04%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.887.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.447"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.20%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.799"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.69%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.36"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +4.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "149.69"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.51%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "132.89"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.49%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.896.57"
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.237"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +26.44%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.820"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +7.41%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.227"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.34%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.38%  "

$ws.Range("E35").Value = "  +4.34%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.08783"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.82%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.611"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.27%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06646"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.12%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "9.205"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02408"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.16%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2231"
$ws.Range("D41").Style = "Normal"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.276"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.46%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.6442"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9990"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.13%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.84"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.40%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6146"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +3.22%  "

$ws.Range("E47").Value = "  +1.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.124"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.68%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.67"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.68%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07315"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.74%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.90"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.98%  "
